$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 464.85715
$ws.Range("J4").Value = 774.75
$ws.Range("L4").Value = 774.75
$ws.Range("N4").Value = -1002.75

$ws.Range("H40").Value = 6129.5
$ws.Range("J40").Value = 7499
$ws.Range("L40").Value = 7499
$ws.Range("N40").Value = -7849

$ws.Range("H88").Value = 3356.25
$ws.Range("J88").Value = 3692.8572
$ws.Range("L88").Value = 3692.8572
$ws.Range("N88").Value = -4504.8572

$ws.Range("H91").Value = 3356.25
$ws.Range("J91").Value = 3692.8572
$ws.Range("L91").Value = 3692.8572
$ws.Range("N91").Value = -6500.8572

$ws.Range("H98").Value = 1012.5294
$ws.Range("J98").Value = 2066.3333
$ws.Range("L98").Value = 2066.3333
$ws.Range("N98").Value = -5062.3333

$ws.Range("H122").Value = 1012.5294
$ws.Range("J122").Value = 2066.3333
$ws.Range("L122").Value = 6198.999899999999
$ws.Range("N122").Value = -11098.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6339.127
$ws.Range("I32").Value = 470.2131
$ws.Range("K32").Value = 470.2131
$ws.Range("M32").Value = -183.2131

$ws.Range("H88").Value = 1156.8
$ws.Range("I88").Value = 825.6
$ws.Range("K88").Value = 825.6
$ws.Range("M88").Value = -419.6

$ws.Range("H91").Value = 1156.8
$ws.Range("I91").Value = 825.6
$ws.Range("K91").Value = 825.6
$ws.Range("M91").Value = 578.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2847
$ws.Range("I94").Value = 1311.3
$ws.Range("J94").Value = 4243.091
$ws.Range("K94").Value = 1311.3
$ws.Range("L94").Value = 4243.091
$ws.Range("M94").Value = -860.3
$ws.Range("N94").Value = -5145.091

$ws.Range("H107").Value = 3631.5
$ws.Range("I107").Value = 3423.1428
$ws.Range("K107").Value = 3423.1428
$ws.Range("M107").Value = -1503.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5933.2046
$ws.Range("I31").Value = 6207.242
$ws.Range("K31").Value = 6207.242
$ws.Range("M31").Value = -5912.242

$ws.Range("H34").Value = 5933.2046
$ws.Range("I34").Value = 6207.242
$ws.Range("K34").Value = 6207.242
$ws.Range("M34").Value = -6005.242

$ws.Range("H58").Value = 1586.875
$ws.Range("I58").Value = 1592.6
$ws.Range("J58").Value = 1577.3334
$ws.Range("K58").Value = 1592.6
$ws.Range("L58").Value = 1577.3334
$ws.Range("M58").Value = -1389.6
$ws.Range("N58").Value = -1983.3334

$ws.Range("H136").Value = 1586.875
$ws.Range("I136").Value = 1592.6
$ws.Range("J136").Value = 1577.3334
$ws.Range("K136").Value = 4777.799999999999
$ws.Range("L136").Value = 4732.0002
$ws.Range("M136").Value = -2227.799999999999
$ws.Range("N136").Value = -9832.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 759.6667
$ws.Range("J5").Value = 747.6667
$ws.Range("L5").Value = 2243.0001
$ws.Range("N5").Value = -2467.0001

$ws.Range("H14").Value = 1444.7858
$ws.Range("I14").Value = 1444.7858
$ws.Range("K14").Value = 4334.357400000001
$ws.Range("M14").Value = -4161.357400000001

$ws.Range("H47").Value = 460.57144
$ws.Range("I47").Value = 460.57144
$ws.Range("K47").Value = 1381.71432
$ws.Range("M47").Value = -950.71432

$ws.Range("H56").Value = 10034.25
$ws.Range("I56").Value = 10034.25
$ws.Range("K56").Value = 10034.25
$ws.Range("M56").Value = -9504.25

$ws.Range("H62").Value = 9414.286
$ws.Range("J62").Value = 9414.286
$ws.Range("L62").Value = 28242.858
$ws.Range("N62").Value = -29614.858

$ws.Range("H65").Value = 9414.286
$ws.Range("J65").Value = 9414.286
$ws.Range("L65").Value = 84728.57399999999
$ws.Range("N65").Value = -91592.57399999999

$ws.Range("H74").Value = 9599.556
$ws.Range("I74").Value = 9597
$ws.Range("J74").Value = 9599.875
$ws.Range("K74").Value = 28791
$ws.Range("L74").Value = 28799.625
$ws.Range("M74").Value = -27730
$ws.Range("N74").Value = -30921.625

$ws.Range("H77").Value = 9599.556
$ws.Range("I77").Value = 9597
$ws.Range("J77").Value = 9599.875
$ws.Range("K77").Value = 86373
$ws.Range("L77").Value = 86398.875
$ws.Range("M77").Value = -81069
$ws.Range("N77").Value = -97006.875

$ws.Range("H113").Value = 1822.5385
$ws.Range("J113").Value = 1895.25
$ws.Range("L113").Value = 5685.75
$ws.Range("N113").Value = -10025.75

$ws.Range("H135").Value = 759.6667
$ws.Range("J135").Value = 747.6667
$ws.Range("L135").Value = 6729.0003
$ws.Range("N135").Value = -11799.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4503.3335
$ws.Range("I19").Value = 3005
$ws.Range("J19").Value = 7500
$ws.Range("K19").Value = 3005
$ws.Range("L19").Value = 7500
$ws.Range("M19").Value = -2717
$ws.Range("N19").Value = -8076

$ws.Range("H80").Value = 3687
$ws.Range("I80").Value = 2798.6
$ws.Range("J80").Value = 10350
$ws.Range("K80").Value = 2798.6
$ws.Range("L80").Value = 10350
$ws.Range("M80").Value = -1800.6
$ws.Range("N80").Value = -12346

$ws.Range("H83").Value = 3687
$ws.Range("I83").Value = 2798.6
$ws.Range("J83").Value = 10350
$ws.Range("K83").Value = 13993
$ws.Range("L83").Value = 51750
$ws.Range("M83").Value = -9001
$ws.Range("N83").Value = -61734

$ws.Range("H102").Value = 2706.8823
$ws.Range("I102").Value = 2719.8125
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2719.8125
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1097.8125
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6199.909
$ws.Range("J100").Value = 7187.5
$ws.Range("L100").Value = 7187.5
$ws.Range("N100").Value = -8269.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3930.3333
$ws.Range("I81").Value = 4195.5
$ws.Range("K81").Value = 8391
$ws.Range("M81").Value = -7330

$ws.Range("H84").Value = 3930.3333
$ws.Range("I84").Value = 4195.5
$ws.Range("K84").Value = 41955
$ws.Range("M84").Value = -36651

$ws.Range("H96").Value = 1436.6666
$ws.Range("I96").Value = 1655
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1655
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = -282
$ws.Range("N96").Value = -3746
